$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: data edits (status/priority changes) ---
# "Implement final GUI design" row (originally Id 8): Status In progress -> Completed
$ws.Range("D9").Value = $ws.Range("D2").Value()
# "Implement themes" row (originally Id 14): Status In progress -> Completed
$ws.Range("D15").Value = $ws.Range("D2").Value()
# "BUG: Clicking on a remove track does not stop playing file" row (Id 26): Priority blank -> Medium
$ws.Range("C27").Value = $ws.Range("C3").Value()

# --- Step 2: convert column A helper formulas (=A(n-1)+1) into static values ---
$rngA = $ws.Range("A2:A30")
$rngA.Value = $rngA.Value()

# --- Step 3: sort ascending by Status to regroup, then re-number Id sequentially ---
$rngAll = $ws.Range("A2:F30")
$keyD = $ws.Range("D2")
$rngAll.Sort($keyD, 1, $null, $null, 1, $null, 1, 2, $false, $null, $null, 1)

for ($i = 2; $i -le 28; $i++) {
    $ws.Cells.Item($i, 1).Value = $i - 1
}

# --- Step 4: sort descending by Status for the final display order ---
$rngAll.Sort($keyD, 2, $null, $null, 1, $null, 1, 2, $false, $null, $null, 1)

Write-Host "done"
